# ContactPage.xlsx update — "Updating repo with latestCodes"
#
# Adds 6 new locator rows (Zip / Province / Privacy-consent fields) to the
# ContactPage locator dictionary and switches the workbook off manual
# calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook previously forced manual recalculation; the update turns
# that off (keeps automatic calculation).
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- New rows -----------------------------------------------------------
# Cell values are written in the exact order the original author entered
# them so that the shared-string table is rebuilt with the same ordering.

# Row 23: PostalCode text box
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = "input#PostalCode"
$ws.Range("B23").Value = "ContactPage_TextBox_Zip"

# Row 24: invalid-zip error message (Text formatted locator name)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "ContactPage_ErrorMessage_InvalidZip_TextBox_Zip"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("C24").Value = ".mktoError #ValidMsgPostalCode"

# Row 25: State/Province dropdown
$ws.Range("A25").Value = 24
$ws.Range("C25").Value = "select#State"
$ws.Range("B25").Value = "ContactPage_Dropdown_Province"

# Row 26: invalid-province error message (Text formatted locator name)
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "ContactPage_ErrorMessage_InvalidProvince_Dropdown_Province"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("C26").Value = ".mktoError #ValidMsgState"

# Row 27: privacy consent checkbox (Text formatted locator name)
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "ContactPage_CheckBox_PrivacyConsent"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("C27").Value = "input#mktoCheckbox_142098_0"

# Row 28: privacy consent label (Text formatted locator name)
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "ContactPage_Label_PrivacyConsent"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("C28").Value = "label#LblmktoCheckbox_142098_0"

# --- View state -----------------------------------------------------------
# Selection moves to C12 and the window is scrolled so row 7 is at the top.
[void]$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
